# Revert "remove double slashes": restore double backslashes in the
# urbansim_path-style directory strings in column L of the "all_runs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L11").Value = "A:\\Projects\\2035_TM160_NGFr2_NP04_Path1_02"
$ws.Range("L14").Value = "G:\\Projects\\2035_TM160_NGF_r2_NoProject_01"
$ws.Range("L19").Value = "B:\\Projects\\2035_TM160_NGF_r2_NoProject_04"
$ws.Range("L22").Value = "A:\\Projects\\2035_TM160_NGFr2_NP04_Path4_01"
$ws.Range("L25").Value = "A:\\Projects\\2035_TM160_NGFr2_NP04_Path4_01"
$ws.Range("L28").Value = "F:\\Projects\\2035_TM160_NGFr2_NP04_Path5_01"
$ws.Range("L31").Value = "H:\\Projects\\2035_TM160_NGFr2_NP04_Path6_01"
$ws.Range("L34").Value = "H:\\Projects\\2035_TM160_NGFr2_NP04_Path6_02"

# Restore the view: move the active selection to A32 instead of L34
# (this also clears the stale "H1" top-left scroll position).
$ws.Range("A32").Select()
